$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.217190146446228
$ws.Range("B1").Value = 0.9544359445571899
$ws.Range("C1").Value = 4.18715238571167
$ws.Range("D1").Value = 2.638900756835938
$ws.Range("E1").Value = 0.770750880241394
